$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 52 (shifts old rows 52:72 down to 53:73)
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row with the De Havilland Comet 4 data
$ws.Range("A52").Value = "De Havilland"
$ws.Range("B52").Value = "Comet 4"
$ws.Range("C52").Value = "Narrow"
$ws.Range("D52").Value = 1949
$ws.Range("G52").Value = 109
$ws.Range("I52").Value = 30.62663043478261
$ws.Range("J52").Value = 0.1846523950527501
$ws.Range("K52").Value = 3.451363636363636
$ws.Range("L52").Value = 313.8715596330275
$ws.Range("M52").Value = 15.01504580607603
